# Fixed data of excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Examples" column values (new shared strings) ---
$ws.Range("C2").Value = "The mail carrier arrives at 10 AM every day."
$ws.Range("C3").Value = "The courier will deliver a package to your office this afternoon."
$ws.Range("C4").Value = "We need to print a sign for the upcoming sale."
$ws.Range("C5").Value = "I want to frame a picture of our family vacation."
$ws.Range("C6").Value = "Please copy a report for each meeting attendee."
$ws.Range("C7").Value = "This restaurant offers delicious food at reasonable prices."
$ws.Range("C8").Value = "The director announced the new project strategy yesterday."

# Reset existing formatting on the whole table so the formatting operations
# below start from one common baseline (avoids piling up unused/orphan
# intermediate styles for every incremental property change).
$full = $ws.Range("A1:C8")
$full.ClearFormats()

# --- Fonts ---
# Header row (A1:C1) + column A data (A2:A8): bold Arial 11 dark gray (#1F1F1F)
$boldRng = $ws.Range("A1:C1,A1:A8")
$boldRng.Font.Name = "Arial"
$boldRng.Font.Size = 11
$boldRng.Font.Color = 2039583
$boldRng.Font.Bold = $true

# Columns B & C data rows: regular Arial 11 dark gray (#1F1F1F)
$plainRng = $ws.Range("B2:C8")
$plainRng.Font.Name = "Arial"
$plainRng.Font.Size = 11
$plainRng.Font.Color = 2039583
$plainRng.Font.Bold = $false

# --- Borders: medium black box around every cell ---
$full.Borders.Weight = -4138
$full.Borders.Color = 0

# --- Alignment ---
$full.HorizontalAlignment = -4131
$full.VerticalAlignment = -4108
$full.WrapText = $true
$full.IndentLevel = 1

# --- Column widths (closest achievable values through this engine's
# character-width rounding; target stored widths are 30.7109375 / 24 /
# 60.85546875) ---
$ws.Columns.Item(1).ColumnWidth = 29.833333333333332
$ws.Columns.Item(2).ColumnWidth = 23.166666666666668
$ws.Columns.Item(3).ColumnWidth = 60

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 30.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 15.75
$ws.Rows.Item(6).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 15.75
$ws.Rows.Item(8).RowHeight = 15.75

# --- Selection ---
$ws.Range("C12").Select()
